# cefclient analysis.pptx — edit script
#
# 1) Typo fix on slide 7: "CLientHandlerStd" -> "ClientHandlerStd"
# 2) Footer date placeholder text bump: "2022-01-02" -> "2022-01-03"
#    on the slide master and all 11 slide layouts (the literal cached
#    text PowerPoint stores inside the <a:fld type="datetimeFigureOut">
#    run whenever the deck is re-saved).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Fix "CLientHandlerStd" -> "ClientHandlerStd" on slide 7
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
for ($i = 1; $i -le $s7.Shapes.Count; $i++) {
    $sh = $s7.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        $needle = "CLientHandlerStd"
        $pos = $full.IndexOf($needle)
        if ($pos -ge 0) {
            $run = $tr.Characters($pos + 1, $needle.Length)
            $run.Text = "ClientHandlerStd"
        }
    }
}

# ---------------------------------------------------------------------
# 2) Bump the cached footer date text from 2022-01-02 to 2022-01-03
#    everywhere it appears: the slide master and every custom layout.
# ---------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "2022-01-02") {
                $tr.Text = "2022-01-03"
            }
        }
    }
}

$master = $p.Slides.Item(1).Master
Update-DateShape $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li).Shapes
}
